$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (the todo item row) is being updated:
#  - B3: task-name placeholder text is edited (the trailing "zz" is removed)
#  - E3: status changed from "완료" (done) to "준비" (not started / preparing)
#  - F3: importance changed from 2 to 3 (kept as a text value, like the
#        existing sibling cell F2, not a numeric one)

$ws.Range("B3").Value = "�Է����ּ���"

$ws.Range("E3").Value = "준비"

# Assign F3 as text "3" (quote-prefixed so it is not auto-converted to a
# number), then copy the plain formatting from F2 on top so the cell keeps
# the workbook's original (unstyled) look instead of picking up a new
# "text-quoted" cell style.
$ws.Range("F3").Value = "'3"
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
